$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "27.482.81"
$ws.Range("E2").Value2 = "  -0.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "1.827.45"
$ws.Range("E3").Value2 = "  -1.81%  "
$ws.Range("E4").Value2 = "  -0.93%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "333.20"
$ws.Range("E5").Value2 = "  -0.42%  "
$ws.Range("E6").Value2 = "  -0.84%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.4565"
$ws.Range("E7").Value2 = "  -1.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.3832"
$ws.Range("E8").Value2 = "  -2.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "46.14"
$ws.Range("E9").Value2 = "  -0.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.07856"
$ws.Range("E10").Value2 = "  -1.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.9580"
$ws.Range("E11").Value2 = "  -4.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "21.05"
$ws.Range("E12").Value2 = "  -2.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "5.837"
$ws.Range("E13").Value2 = "  -1.88%  "
$ws.Range("B14").Value2 = "Chainlink"
$ws.Range("C14").Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "7.051"
$ws.Range("E14").Value2 = "  -2.23%  "
$ws.Range("B15").Value2 = "WrappedEther"
$ws.Range("C15").Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "1.796.67"
$ws.Range("E15").Value2 = "  -3.85%  "
$ws.Range("E16").Value2 = "  -0.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "89.48"
$ws.Range("E17").Value2 = "  +1.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "0.06567"
$ws.Range("E18").Value2 = "  -2.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "0.00001022"
$ws.Range("E19").Value2 = "  -1.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "17.10"
$ws.Range("E20").Value2 = "  -0.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "1.003"
$ws.Range("E21").Value2 = "  -0.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "27.457.89"
$ws.Range("E22").Value2 = "  -0.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "5.286"
$ws.Range("E23").Value2 = "  -2.85%  "
$ws.Range("E24").Value2 = "  -1.40%  "
$ws.Range("E25").Value2 = "  -2.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "159.30"
$ws.Range("E26").Value2 = "  -0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "2.025.71"
$ws.Range("E27").Value2 = "  -2.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "19.37"
$ws.Range("E28").Value2 = "  -1.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "2.042"
$ws.Range("E29").Value2 = "  -4.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "5.271"
$ws.Range("E30").Value2 = "  -3.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "117.90"
$ws.Range("E31").Value2 = "  -3.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "0.09358"
$ws.Range("E32").Value2 = "  -0.40%  "
$ws.Range("E33").Value2 = "  -4.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "3.572"
$ws.Range("E34").Value2 = "  -1.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "5.221"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "1.317"
$ws.Range("E36").Value2 = "  -1.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "0.05900"
$ws.Range("E37").Value2 = "  -1.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "0.02186"
$ws.Range("E38").Value2 = "  -2.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "8.081"
$ws.Range("E39").Value2 = "  -2.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "1.003"
$ws.Range("E40").Value2 = "  -0.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "1.144"
$ws.Range("E41").Value2 = "  -4.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "0.5731"
$ws.Range("E42").Value2 = "  -3.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.1817"
$ws.Range("E43").Value2 = "  -2.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "9.899"
$ws.Range("E44").Value2 = "  -4.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "1.263"
$ws.Range("E45").Value2 = "  +1.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.5397"
$ws.Range("E46").Value2 = "  -3.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "11.72"
$ws.Range("E47").Value2 = "  -3.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "1.885"
$ws.Range("E48").Value2 = "  -1.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.06929"
$ws.Range("E49").Value2 = "  +2.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "110.29"
$ws.Range("E50").Value2 = "  -1.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "1.003"
$ws.Range("E51").Value2 = "  -33.19%  "
